$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 222, shifting rows 222:276 down to 223:277
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with new data
$ws.Range("A222").Value = 4
$ws.Range("B222").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C222").Value = "Los Lagos"
$ws.Range("D222").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D222").Value = 44722
$ws.Range("E222").Value = 10
$ws.Range("F222").Value = 100112037
$ws.Range("G222").Value = "Cebollín"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 150
$ws.Range("K222").Value = 10000
$ws.Range("L222").Value = 10000
$ws.Range("M222").Value = 10000
$ws.Range("N222").Value = "`$/paquete 36 unidades"
$ws.Range("O222").Value = "Región Metropolitana"
$ws.Range("P222").Value = 278
$ws.Range("Q222").Value = 36
$ws.Range("R222").Value = "Hortaliza"
